$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the policy number (NroPoliza) in E2
$ws.Range("E2").Value = 12112002082

# Update the claim date (FechaSiniestro) in G2 (stored as quote-prefixed text,
# leading apostrophe keeps it as text rather than being parsed as a date)
$ws.Range("G2").Value = "'19/03/2021"

# Select E2 (view scrolls back so the sheet starts at column A again)
$ws.Range("E2").Select()
